$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update header row (row 1) with corrected/reordered column titles
$ws.Range("A1").Value = "Corriente (mA)"
$ws.Range("B1").Value = "Voltaje Hall a 1,5 A (mV)"
$ws.Range("C1").Value = "Voltaje Hall a 1,75 A (mV)"
$ws.Range("D1").Value = "Voltaje Hall a 2 A (mV)"
$ws.Range("E1").Value = "Voltaje Hall a 2,25 A (mV)"
$ws.Range("F1").Value = "Voltaje Hall a 2,5 A (mV)"

# Reorder data columns (B..F) so values are sorted ascending by field current (1.5A..2.5A)
$ws.Range("B2").Value = 54.3
$ws.Range("C2").Value = 59.4
$ws.Range("D2").Value = 70.400000000000006
$ws.Range("E2").Value = 69.2
$ws.Range("F2").Value = 73.7
$ws.Range("B3").Value = 49.8
$ws.Range("C3").Value = 53.4
$ws.Range("D3").Value = 62
$ws.Range("E3").Value = 65.3
$ws.Range("F3").Value = 66
$ws.Range("B4").Value = 45
$ws.Range("C4").Value = 51.8
$ws.Range("D4").Value = 58.1
$ws.Range("E4").Value = 56.8
$ws.Range("F4").Value = 60.8
$ws.Range("B5").Value = 40.799999999999997
$ws.Range("C5").Value = 43.3
$ws.Range("D5").Value = 50.8
$ws.Range("E5").Value = 51
$ws.Range("F5").Value = 55.3
$ws.Range("B6").Value = 33.700000000000003
$ws.Range("C6").Value = 38.299999999999997
$ws.Range("D6").Value = 44.7
$ws.Range("E6").Value = 45
$ws.Range("F6").Value = 48.8
$ws.Range("B7").Value = 28.8
$ws.Range("C7").Value = 34.1
$ws.Range("D7").Value = 38.700000000000003
$ws.Range("E7").Value = 40.299999999999997
$ws.Range("F7").Value = 42.9
$ws.Range("B8").Value = 23.7
$ws.Range("C8").Value = 26.5
$ws.Range("D8").Value = 33.6
$ws.Range("E8").Value = 33.200000000000003
$ws.Range("F8").Value = 36.1
$ws.Range("B9").Value = 19.600000000000001
$ws.Range("C9").Value = 22.6
$ws.Range("D9").Value = 28.1
$ws.Range("E9").Value = 26.2
$ws.Range("F9").Value = 31.7
$ws.Range("B10").Value = 14.2
$ws.Range("C10").Value = 17.399999999999999
$ws.Range("D10").Value = 22.2
$ws.Range("E10").Value = 19.399999999999999
$ws.Range("F10").Value = 25
$ws.Range("B11").Value = 9.4
$ws.Range("C11").Value = 11.4
$ws.Range("D11").Value = 16
$ws.Range("E11").Value = 14.8
$ws.Range("F11").Value = 18.5
$ws.Range("B12").Value = 4.4000000000000004
$ws.Range("C12").Value = 7.5
$ws.Range("D12").Value = 9.1999999999999993
$ws.Range("E12").Value = 8.3000000000000007
$ws.Range("F12").Value = 11.1
$ws.Range("B13").Value = -0.8
$ws.Range("C13").Value = -0.7
$ws.Range("D13").Value = 2.8
$ws.Range("E13").Value = 4.8
$ws.Range("F13").Value = 4.7
$ws.Range("B14").Value = -5.2
$ws.Range("C14").Value = -4.8
$ws.Range("D14").Value = -3.8
$ws.Range("E14").Value = -2.9
$ws.Range("F14").Value = -0.4
$ws.Range("B15").Value = -11.2
$ws.Range("C15").Value = -11.3
$ws.Range("D15").Value = -9.1
$ws.Range("E15").Value = -10
$ws.Range("F15").Value = -6.9
$ws.Range("B16").Value = -17.2
$ws.Range("C16").Value = -18
$ws.Range("D16").Value = -15.2
$ws.Range("E16").Value = -16.8
$ws.Range("F16").Value = -15.1
$ws.Range("B17").Value = -19.7
$ws.Range("C17").Value = -24.5
$ws.Range("D17").Value = -20.2
$ws.Range("E17").Value = -25.8
$ws.Range("F17").Value = -20.9
$ws.Range("B18").Value = -28
$ws.Range("C18").Value = -32.200000000000003
$ws.Range("D18").Value = -29.5
$ws.Range("E18").Value = -33.700000000000003
$ws.Range("F18").Value = -29.7
$ws.Range("B19").Value = -34.299999999999997
$ws.Range("C19").Value = -36.700000000000003
$ws.Range("D19").Value = -36.6
$ws.Range("E19").Value = -39.700000000000003
$ws.Range("F19").Value = -38.799999999999997
$ws.Range("B20").Value = -40.4
$ws.Range("C20").Value = -42.8
$ws.Range("D20").Value = -41.4
$ws.Range("E20").Value = -46.3
$ws.Range("F20").Value = -43.2
$ws.Range("B21").Value = -46.4
$ws.Range("C21").Value = -48.5
$ws.Range("D21").Value = -46
$ws.Range("E21").Value = -50.7
$ws.Range("F21").Value = -52.9
$ws.Range("B22").Value = -49.8
$ws.Range("C22").Value = -52.9
$ws.Range("D22").Value = -54
$ws.Range("E22").Value = -57.3
$ws.Range("F22").Value = -57.8
$ws.Range("B23").Value = -53.5
$ws.Range("C23").Value = -57
$ws.Range("D23").Value = -60.8
$ws.Range("E23").Value = -62.6
$ws.Range("F23").Value = -63.9
$ws.Range("B24").Value = -60.4
$ws.Range("C24").Value = -63.3
$ws.Range("D24").Value = -64.900000000000006
$ws.Range("E24").Value = -68.3
$ws.Range("F24").Value = -68.900000000000006
$ws.Range("B25").Value = -64.599999999999994
$ws.Range("C25").Value = -70
$ws.Range("D25").Value = -72.400000000000006
$ws.Range("E25").Value = -76.3
$ws.Range("F25").Value = -74.5
$ws.Range("B26").Value = -70
$ws.Range("C26").Value = -76.099999999999994
$ws.Range("D26").Value = -75.8
$ws.Range("E26").Value = -82.4
$ws.Range("F26").Value = -82.3
$ws.Range("B27").Value = -75.900000000000006
$ws.Range("C27").Value = -80.400000000000006
$ws.Range("D27").Value = -86.7
$ws.Range("E27").Value = -90.2
$ws.Range("F27").Value = -86.9
$ws.Range("B28").Value = -79.599999999999994
$ws.Range("C28").Value = -86.5
$ws.Range("D28").Value = -92.6
$ws.Range("E28").Value = -92.1
$ws.Range("F28").Value = -95.5
$ws.Range("B29").Value = -86.1
$ws.Range("C29").Value = -92.3
$ws.Range("D29").Value = -96.2
$ws.Range("E29").Value = -100.4
$ws.Range("F29").Value = -101.4
$ws.Range("B30").Value = -92.3
$ws.Range("C30").Value = -98.4
$ws.Range("D30").Value = -102.4
$ws.Range("E30").Value = -107.6
$ws.Range("F30").Value = -106.9
$ws.Range("B31").Value = -99.1
$ws.Range("C31").Value = -100.9
$ws.Range("D31").Value = -108.3
$ws.Range("E31").Value = -114.9
$ws.Range("F31").Value = -112.9
$ws.Range("B32").Value = -104.9
$ws.Range("C32").Value = -107.2
$ws.Range("D32").Value = -115.6
$ws.Range("E32").Value = -120.4
$ws.Range("F32").Value = -122.2
